# Selenabler AMS test data workbook update (commit: "Last commit on 21Apr2023")
#
# Summary of the change being applied:
#  - Add a new "User" worksheet (as the last tab) with login/user test data
#  - Update the remembered cell-selection on several existing sheets
#  - Widen column B on the "EmployeeModule" sheet
#
# NOTE on shared-string ordering: the target workbook's sharedStrings.xml has
# the "User" sheet's column-B values (the data) inserted before its column-A
# values (the labels). The engine appends new shared strings strictly in the
# order cell values are assigned, so we intentionally write column B first
# (in a specific row order), then column A, to reproduce that exact ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. LoginSheet: move the remembered selection from E14 to A15
# ---------------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("LoginSheet")
$loginSheet.Range("A15").Select()

# ---------------------------------------------------------------------------
# 2. AssetModule: move the remembered selection from A14 to A10
# ---------------------------------------------------------------------------
$assetModule = $wb.Worksheets.Item("AssetModule")
$assetModule.Range("A10").Select()

# ---------------------------------------------------------------------------
# 3. EmployeeModule: widen column B (~25.43 -> 42 chars) and move the
#    remembered selection from E13 to E7
# ---------------------------------------------------------------------------
$employeeModule = $wb.Worksheets.Item("EmployeeModule")
$employeeModule.Columns.Item(2).ColumnWidth = 41.15
$employeeModule.Range("E7").Select()

# ---------------------------------------------------------------------------
# 4. AssetVendor: move the remembered selection to E17
# ---------------------------------------------------------------------------
$assetVendor = $wb.Worksheets.Item("AssetVendor")
$assetVendor.Range("E17").Select()

# ---------------------------------------------------------------------------
# 5. AssetsStatus: move the remembered selection from D7 to C12
#    (this sheet was previously the active tab; it no longer is once the
#    new "User" sheet below is created/activated)
# ---------------------------------------------------------------------------
$assetsStatus = $wb.Worksheets.Item("AssetsStatus")
$assetsStatus.Range("C12").Select()

# ---------------------------------------------------------------------------
# 6. Add the new "User" sheet as the last tab and populate it
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$userSheet = $wb.Worksheets.Add($null, $lastSheet)
$userSheet.Name = "User"

# Column B (values) first, in this specific row order, then column A
# (labels) -- see note above about shared-string ordering.
$userSheet.Range("B1").Value = "P1345365"
$userSheet.Range("B2").Value = "AkashD"
$userSheet.Range("B7").Value = "akash.dhole@in.ncs-i.com"
$userSheet.Range("B8").Value = "USER"
$userSheet.Range("B5").Value = "akash@123"
$userSheet.Range("B6").Value = "AkashDhole"
$userSheet.Range("B3").Value = "akash.dhole@ncs.com.sg"
$userSheet.Range("B4").Value = "ADMIN"

$userSheet.Range("A1").Value = "EmpID1"
$userSheet.Range("A2").Value = "UserName1"
$userSheet.Range("A3").Value = "Email1"
$userSheet.Range("A4").Value = "Role1"
$userSheet.Range("A5").Value = "UserPassword1"
$userSheet.Range("A6").Value = "NewUsername1"
$userSheet.Range("A7").Value = "NewEmail1"
$userSheet.Range("A8").Value = "ChangeRole1"

$userSheet.Range("A9").Value = "EOF"
$userSheet.Range("B9").Value = "EOF"

# Column B uses the "Text" number format (same style already used by the
# other sheets' value columns).
$userSheet.Range("B1:B9").NumberFormat = "@"

# Column widths roughly matching the source sheet (closest attainable via
# the character-width COM API).
$userSheet.Columns.Item(1).ColumnWidth = 17.25
$userSheet.Columns.Item(2).ColumnWidth = 26.59

$userSheet.PageSetup.Orientation = 1

# Finally, set the remembered selection on the new sheet -- doing this last
# makes "User" the active tab, matching the target workbook.
$userSheet.Range("E14").Select()
